$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 259 (shifts old rows 259..312 down to 260..313).
$ws.Rows.Item(259).Insert()

# Populate the newly inserted row 259 with the new weekly price observation.
$ws.Range("A259").Value = 4
$ws.Range("B259").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C259").Value = "Los Lagos"
$ws.Range("D259").Value = 44798
$ws.Range("E259").Value = 10
$ws.Range("F259").Value = 100112043
$ws.Range("G259").Value = "Pepino ensalada"
$ws.Range("H259").Value = "Sin especificar"
$ws.Range("I259").Value = "Primera"
$ws.Range("J259").Value = 150
$ws.Range("K259").Value = 27000
$ws.Range("L259").Value = 27000
$ws.Range("M259").Value = 27000
$ws.Range("N259").Value = "`$/caja 60 unidades"
$ws.Range("O259").Value = "Región de Arica y Parinacota"
$ws.Range("P259").Value = 450
$ws.Range("Q259").Value = 60
$ws.Range("R259").Value = "Hortaliza"

# Match the date-format style used by the rest of column D.
$ws.Range("D259").NumberFormat = $ws.Range("D260").NumberFormat()
